$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.472070813179016
$ws.Range("B1").Value = 1.738719463348389
$ws.Range("C1").Value = 1.864494204521179
$ws.Range("D1").Value = 2.128275871276855
$ws.Range("E1").Value = 2.693896055221558
